$wb = $excel.ActiveWorkbook

# --- Sheet1: Daily Orders ---
# A new order (#21, Pooja) came in after the existing order (#20, Minakshi).
# Orders are listed most-recent-first, so insert a fresh row 2 and push the
# existing row 2 down to row 3.
$ws1 = $wb.Worksheets.Item("Daily Orders")
$ws1.Rows.Item(2).Insert()

$ws1.Cells.Item(2, 1).Value = 21
$ws1.Cells.Item(2, 2).Value = "2026-01-20 05:33"
$ws1.Cells.Item(2, 3).Value = "Pooja"
$ws1.Cells.Item(2, 4).Value = "A 1608"
$ws1.Cells.Item(2, 5).NumberFormat = "@"
$ws1.Cells.Item(2, 5).Value = "9096648553"
$ws1.Cells.Item(2, 6).Value = "Jawar Bhakari x1"
$ws1.Cells.Item(2, 7).Value = 20
$ws1.Cells.Item(2, 8).Value = "NEW"
$ws1.Cells.Item(2, 9).Value = "PENDING"
$ws1.Cells.Item(2, 10).NumberFormat = "@"
$ws1.Cells.Item(2, 10).Value = "2026-01-20"
$ws1.Cells.Item(2, 11).Value = "11:03"

# --- Sheet2: Summary ---
# One more order overall, one more in the "New" bucket, and the new
# order's total ($20) adds to total revenue (paid amount is unaffected
# since the new order is still PENDING payment).
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Cells.Item(2, 1).Value = 2
$ws2.Cells.Item(2, 2).Value = 1
$ws2.Cells.Item(2, 7).Value = 50

# --- Sheet3: Items Breakdown ---
# New line item for the Jawar Bhakari ordered in the new order.
$ws3 = $wb.Worksheets.Item("Items Breakdown")
$ws3.Cells.Item(3, 1).Value = "Jawar Bhakari"
$ws3.Cells.Item(3, 2).Value = 1
$ws3.Cells.Item(3, 3).Value = 20
